$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update/extend the Ajo price-report rows (339-350) per the latest weekly
# fruit/vegetable data drop. Rows 339-345 are rewritten in place and rows
# 346-350 are appended, which naturally grows the sheet dimension to R350.

# Row 339
$ws.Cells.Item(339, 1).Value = 10
$ws.Cells.Item(339, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(339, 3).Value = 'La Araucanía'
$ws.Cells.Item(339, 4).Value = 44448
$ws.Cells.Item(339, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(339, 5).Value = 9
$ws.Cells.Item(339, 6).Value = 100112003
$ws.Cells.Item(339, 7).Value = 'Ajo'
$ws.Cells.Item(339, 8).Value = 'Chino'
$ws.Cells.Item(339, 9).Value = 'Primera'
$ws.Cells.Item(339, 10).Value = 255
$ws.Cells.Item(339, 11).Value = 17000
$ws.Cells.Item(339, 12).Value = 17000
$ws.Cells.Item(339, 13).Value = 17000
$ws.Cells.Item(339, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(339, 15).Value = 'China'
$ws.Cells.Item(339, 16).Value = 1700
$ws.Cells.Item(339, 17).Value = 10
$ws.Cells.Item(339, 18).Value = 'Hortaliza'

# Row 340
$ws.Cells.Item(340, 1).Value = 10
$ws.Cells.Item(340, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(340, 3).Value = 'La Araucanía'
$ws.Cells.Item(340, 4).Value = 44448
$ws.Cells.Item(340, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(340, 5).Value = 9
$ws.Cells.Item(340, 6).Value = 100112003
$ws.Cells.Item(340, 7).Value = 'Ajo'
$ws.Cells.Item(340, 8).Value = 'Chino'
$ws.Cells.Item(340, 9).Value = 'Primera'
$ws.Cells.Item(340, 10).Value = 285
$ws.Cells.Item(340, 11).Value = 20000
$ws.Cells.Item(340, 12).Value = 20000
$ws.Cells.Item(340, 13).Value = 20000
$ws.Cells.Item(340, 14).Value = '$/malla 10 kilos'
$ws.Cells.Item(340, 15).Value = 'China'
$ws.Cells.Item(340, 16).Value = 2000
$ws.Cells.Item(340, 17).Value = 10
$ws.Cells.Item(340, 18).Value = 'Hortaliza'

# Row 341
$ws.Cells.Item(341, 1).Value = 10
$ws.Cells.Item(341, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(341, 3).Value = 'La Araucanía'
$ws.Cells.Item(341, 4).Value = 44448
$ws.Cells.Item(341, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(341, 5).Value = 9
$ws.Cells.Item(341, 6).Value = 100112003
$ws.Cells.Item(341, 7).Value = 'Ajo'
$ws.Cells.Item(341, 8).Value = 'Rosado'
$ws.Cells.Item(341, 9).Value = '1a (guarda)'
$ws.Cells.Item(341, 10).Value = 155
$ws.Cells.Item(341, 11).Value = 10000
$ws.Cells.Item(341, 12).Value = 10000
$ws.Cells.Item(341, 13).Value = 10000
$ws.Cells.Item(341, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(341, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(341, 16).Value = 1000
$ws.Cells.Item(341, 17).Value = 10
$ws.Cells.Item(341, 18).Value = 'Hortaliza'

# Row 342
$ws.Cells.Item(342, 1).Value = 10
$ws.Cells.Item(342, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(342, 3).Value = 'La Araucanía'
$ws.Cells.Item(342, 4).Value = 44448
$ws.Cells.Item(342, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(342, 5).Value = 9
$ws.Cells.Item(342, 6).Value = 100112003
$ws.Cells.Item(342, 7).Value = 'Ajo'
$ws.Cells.Item(342, 8).Value = 'Rosado'
$ws.Cells.Item(342, 9).Value = '2a (guarda)'
$ws.Cells.Item(342, 10).Value = 45
$ws.Cells.Item(342, 11).Value = 3000
$ws.Cells.Item(342, 12).Value = 3000
$ws.Cells.Item(342, 13).Value = 3000
$ws.Cells.Item(342, 14).Value = '$/trenza 50 unidades'
$ws.Cells.Item(342, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(342, 16).Value = 600
$ws.Cells.Item(342, 17).Value = 5
$ws.Cells.Item(342, 18).Value = 'Hortaliza'

# Row 343
$ws.Cells.Item(343, 1).Value = 10
$ws.Cells.Item(343, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(343, 3).Value = 'La Araucanía'
$ws.Cells.Item(343, 4).Value = 44448
$ws.Cells.Item(343, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(343, 5).Value = 9
$ws.Cells.Item(343, 6).Value = 100112003
$ws.Cells.Item(343, 7).Value = 'Ajo'
$ws.Cells.Item(343, 8).Value = 'Rosado'
$ws.Cells.Item(343, 9).Value = '3a (guarda)'
$ws.Cells.Item(343, 10).Value = 35
$ws.Cells.Item(343, 11).Value = 2500
$ws.Cells.Item(343, 12).Value = 2500
$ws.Cells.Item(343, 13).Value = 2500
$ws.Cells.Item(343, 14).Value = '$/trenza 50 unidades'
$ws.Cells.Item(343, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(343, 16).Value = 500
$ws.Cells.Item(343, 17).Value = 5
$ws.Cells.Item(343, 18).Value = 'Hortaliza'

# Row 344
$ws.Cells.Item(344, 1).Value = 10
$ws.Cells.Item(344, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(344, 3).Value = 'La Araucanía'
$ws.Cells.Item(344, 4).Value = 44167
$ws.Cells.Item(344, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(344, 5).Value = 9
$ws.Cells.Item(344, 6).Value = 100112003
$ws.Cells.Item(344, 7).Value = 'Ajo'
$ws.Cells.Item(344, 8).Value = 'Chino'
$ws.Cells.Item(344, 9).Value = 'Primera'
$ws.Cells.Item(344, 10).Value = 175
$ws.Cells.Item(344, 11).Value = 9000
$ws.Cells.Item(344, 12).Value = 10000
$ws.Cells.Item(344, 13).Value = 9371
$ws.Cells.Item(344, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(344, 15).Value = 'China'
$ws.Cells.Item(344, 16).Value = 937
$ws.Cells.Item(344, 17).Value = 10
$ws.Cells.Item(344, 18).Value = 'Hortaliza'

# Row 345
$ws.Cells.Item(345, 1).Value = 10
$ws.Cells.Item(345, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(345, 3).Value = 'La Araucanía'
$ws.Cells.Item(345, 4).Value = 44238
$ws.Cells.Item(345, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(345, 5).Value = 9
$ws.Cells.Item(345, 6).Value = 100112003
$ws.Cells.Item(345, 7).Value = 'Ajo'
$ws.Cells.Item(345, 8).Value = 'Chino'
$ws.Cells.Item(345, 9).Value = 'Primera'
$ws.Cells.Item(345, 10).Value = 450
$ws.Cells.Item(345, 11).Value = 12000
$ws.Cells.Item(345, 12).Value = 13000
$ws.Cells.Item(345, 13).Value = 12444
$ws.Cells.Item(345, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(345, 15).Value = 'China'
$ws.Cells.Item(345, 16).Value = 1244
$ws.Cells.Item(345, 17).Value = 10
$ws.Cells.Item(345, 18).Value = 'Hortaliza'

# Row 346
$ws.Cells.Item(346, 1).Value = 10
$ws.Cells.Item(346, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(346, 3).Value = 'La Araucanía'
$ws.Cells.Item(346, 4).Value = 44399
$ws.Cells.Item(346, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(346, 5).Value = 9
$ws.Cells.Item(346, 6).Value = 100112003
$ws.Cells.Item(346, 7).Value = 'Ajo'
$ws.Cells.Item(346, 8).Value = 'Chino'
$ws.Cells.Item(346, 9).Value = 'Primera'
$ws.Cells.Item(346, 10).Value = 500
$ws.Cells.Item(346, 11).Value = 14000
$ws.Cells.Item(346, 12).Value = 15000
$ws.Cells.Item(346, 13).Value = 14400
$ws.Cells.Item(346, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(346, 15).Value = 'China'
$ws.Cells.Item(346, 16).Value = 1440
$ws.Cells.Item(346, 17).Value = 10
$ws.Cells.Item(346, 18).Value = 'Hortaliza'

# Row 347
$ws.Cells.Item(347, 1).Value = 10
$ws.Cells.Item(347, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(347, 3).Value = 'La Araucanía'
$ws.Cells.Item(347, 4).Value = 44399
$ws.Cells.Item(347, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(347, 5).Value = 9
$ws.Cells.Item(347, 6).Value = 100112003
$ws.Cells.Item(347, 7).Value = 'Ajo'
$ws.Cells.Item(347, 8).Value = 'Chino'
$ws.Cells.Item(347, 9).Value = 'Primera'
$ws.Cells.Item(347, 10).Value = 200
$ws.Cells.Item(347, 11).Value = 18000
$ws.Cells.Item(347, 12).Value = 18000
$ws.Cells.Item(347, 13).Value = 18000
$ws.Cells.Item(347, 14).Value = '$/malla 10 kilos'
$ws.Cells.Item(347, 15).Value = 'China'
$ws.Cells.Item(347, 16).Value = 1800
$ws.Cells.Item(347, 17).Value = 10
$ws.Cells.Item(347, 18).Value = 'Hortaliza'

# Row 348
$ws.Cells.Item(348, 1).Value = 10
$ws.Cells.Item(348, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(348, 3).Value = 'La Araucanía'
$ws.Cells.Item(348, 4).Value = 44399
$ws.Cells.Item(348, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(348, 5).Value = 9
$ws.Cells.Item(348, 6).Value = 100112003
$ws.Cells.Item(348, 7).Value = 'Ajo'
$ws.Cells.Item(348, 8).Value = 'Rosado'
$ws.Cells.Item(348, 9).Value = '1a (guarda)'
$ws.Cells.Item(348, 10).Value = 200
$ws.Cells.Item(348, 11).Value = 14000
$ws.Cells.Item(348, 12).Value = 14000
$ws.Cells.Item(348, 13).Value = 14000
$ws.Cells.Item(348, 14).Value = '$/malla 10 kilos'
$ws.Cells.Item(348, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(348, 16).Value = 1400
$ws.Cells.Item(348, 17).Value = 10
$ws.Cells.Item(348, 18).Value = 'Hortaliza'

# Row 349
$ws.Cells.Item(349, 1).Value = 10
$ws.Cells.Item(349, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(349, 3).Value = 'La Araucanía'
$ws.Cells.Item(349, 4).Value = 44400
$ws.Cells.Item(349, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(349, 5).Value = 9
$ws.Cells.Item(349, 6).Value = 100112003
$ws.Cells.Item(349, 7).Value = 'Ajo'
$ws.Cells.Item(349, 8).Value = 'Chino'
$ws.Cells.Item(349, 9).Value = 'Primera'
$ws.Cells.Item(349, 10).Value = 330
$ws.Cells.Item(349, 11).Value = 14000
$ws.Cells.Item(349, 12).Value = 15000
$ws.Cells.Item(349, 13).Value = 14455
$ws.Cells.Item(349, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(349, 15).Value = 'China'
$ws.Cells.Item(349, 16).Value = 1446
$ws.Cells.Item(349, 17).Value = 10
$ws.Cells.Item(349, 18).Value = 'Hortaliza'

# Row 350
$ws.Cells.Item(350, 1).Value = 10
$ws.Cells.Item(350, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(350, 3).Value = 'La Araucanía'
$ws.Cells.Item(350, 4).Value = 44400
$ws.Cells.Item(350, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(350, 5).Value = 9
$ws.Cells.Item(350, 6).Value = 100112003
$ws.Cells.Item(350, 7).Value = 'Ajo'
$ws.Cells.Item(350, 8).Value = 'Chino'
$ws.Cells.Item(350, 9).Value = 'Primera'
$ws.Cells.Item(350, 10).Value = 200
$ws.Cells.Item(350, 11).Value = 18000
$ws.Cells.Item(350, 12).Value = 18000
$ws.Cells.Item(350, 13).Value = 18000
$ws.Cells.Item(350, 14).Value = '$/malla 10 kilos'
$ws.Cells.Item(350, 15).Value = 'China'
$ws.Cells.Item(350, 16).Value = 1800
$ws.Cells.Item(350, 17).Value = 10
$ws.Cells.Item(350, 18).Value = 'Hortaliza'

